$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.308.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -3.55%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.159.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -8.84%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.04%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'563.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -3.83%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'168.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -5.38%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.613"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -2.26%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.01%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'3.153.76"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -8.93%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -7.09%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -6.50%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.393"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -6.38%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'3.696.06"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -9.05%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  +0.75%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'27.04"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -10.14%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'64.259.48"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -3.38%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  -6.30%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.150.22"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -9.38%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'5.72"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -4.42%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'12.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -7.44%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'351.84"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -5.40%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'7.20"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -6.25%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.31%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'68.04"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -7.50%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.499"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -7.21%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0000116"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -8.52%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'9.60"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -4.05%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -1.63%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -0.05%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -0.15%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -6.16%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'5.45"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -8.24%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'21.89"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -7.65%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'6.59"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -7.09%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -6.74%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("B36").Value = "'ImmutableX"
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").Value = "'1.42"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -9.83%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("B37").Value = "'Monero"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'153.77"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -5.27%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.816"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -7.96%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'26.31"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -6.16%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -6.92%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'2.45"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -5.02%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'2.623.07"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -5.42%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'4.16"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -8.20%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'39.37"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.79%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -8.06%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.0646"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -7.27%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'23.74"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -7.10%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'318.85"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -6.49%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -6.59%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -3.02%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.998"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.06%  "
$ws.Range("E51").Style = "Normal"
